$d = $word.ActiveDocument

# Position at the very end of the document (after "/getid.")
$r = $d.Range($d.Content.End, $d.Content.End)
$r.Collapse(0)

# Insert a blank RTL paragraph followed by a new RTL paragraph containing
# "/group_stats", matching the target OOXML exactly (no direct run
# formatting - just the RTL paragraph mark, same as the existing paragraph).
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr></w:p>' + `
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:t>/group_stats</w:t></w:r></w:p>'

$r.InsertXML($xml)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
